$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '54.288.12'
$ws.Cells.Item(2, 5).Value = '  +0.47%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.268.94'
$ws.Cells.Item(3, 5).Value = '  +0.95%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.18%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '497.82'
$ws.Cells.Item(5, 5).Value = '  +0.75%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '129.21'
$ws.Cells.Item(6, 5).Value = '  +1.20%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.998'
$ws.Cells.Item(7, 5).Value = '  +0.25%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.526'
$ws.Cells.Item(8, 5).Value = '  -0.21%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.0954'
$ws.Cells.Item(9, 5).Value = '  +0.11%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +0.85%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.336'
$ws.Cells.Item(11, 5).Value = '  +3.52%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '4.92'
$ws.Cells.Item(12, 5).Value = '  +5.84%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '23.04'
$ws.Cells.Item(13, 5).Value = '  +6.09%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '2.673.88'
$ws.Cells.Item(14, 5).Value = '  +0.92%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '54.308.12'
$ws.Cells.Item(15, 5).Value = '  +0.66%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +0.68%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '2.270.39'
$ws.Cells.Item(17, 5).Value = '  +0.71%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '10.27'
$ws.Cells.Item(18, 5).Value = '  +2.60%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.15'
$ws.Cells.Item(19, 5).Value = '  +1.36%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '304.05'
$ws.Cells.Item(20, 5).Value = '  +1.45%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.33'
$ws.Cells.Item(21, 5).Value = '  -1.50%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.14%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '60.66'
$ws.Cells.Item(23, 5).Value = '  -2.01%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.997'
$ws.Cells.Item(24, 5).Value = '  -2.12%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.95%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '7.35'
$ws.Cells.Item(26, 5).Value = '  +4.47%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '175.17'

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.0₃0705'
$ws.Cells.Item(28, 5).Value = '  +2.99%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'PancakeSwap'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.61'
$ws.Cells.Item(29, 5).Value = '  +0.45%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'Aptos'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '6.00'
$ws.Cells.Item(30, 5).Value = '  +2.73%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.08'
$ws.Cells.Item(31, 5).Value = '  +1.54%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +0.04%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '17.83'
$ws.Cells.Item(33, 5).Value = '  +1.12%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.998'
$ws.Cells.Item(34, 5).Value = '  +0.32%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.950'
$ws.Cells.Item(35, 5).Value = '  +5.10%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +2.12%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '3.72'
$ws.Cells.Item(37, 5).Value = '  +1.28%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.376'
$ws.Cells.Item(38, 5).Value = '  +1.37%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.40'
$ws.Cells.Item(39, 5).Value = '  +0.51%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Filecoin'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.38'
$ws.Cells.Item(40, 5).Value = '  +0.91%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'RenderToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '4.92'
$ws.Cells.Item(41, 5).Value = '  -0.25%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '125.21'
$ws.Cells.Item(42, 5).Value = '  -0.16%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.0492'
$ws.Cells.Item(43, 5).Value = '  +2.01%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0897'
$ws.Cells.Item(44, 5).Value = '  +1.15%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '244.40'
$ws.Cells.Item(45, 5).Value = '  +3.34%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.547'
$ws.Cells.Item(46, 5).Value = '  +1.08%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.375'
$ws.Cells.Item(47, 5).Value = '  +1.40%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0205'
$ws.Cells.Item(48, 5).Value = '  +1.62%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +0.82%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '16.22'
$ws.Cells.Item(50, 5).Value = '  +0.59%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.52'
$ws.Cells.Item(51, 5).Value = '  +2.28%  '
